$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (avoids Excel auto-coercing numeric-looking
# strings like "596.17" or "1.00" into Double values), while keeping the
# cells style index unchanged (matches original unstyled data cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "67.181.98"
$ws.Range("E2").Value = "  -4.83%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.256.47"
$ws.Range("E3").Value = "  -7.37%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "596.17"
$ws.Range("E5").Value = "  -4.37%  "

# Row 6
Set-TextValue $ws.Range("D6") "150.75"
$ws.Range("E6").Value = "  -12.81%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.249.05"
$ws.Range("E8").Value = "  -7.48%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.540"
$ws.Range("E9").Value = "  -11.52%  "

# Row 10
$ws.Range("E10").Value = "  -14.10%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.68"
$ws.Range("E11").Value = "  -6.76%  "

# Row 12
$ws.Range("E12").Value = "  -14.19%  "

# Row 13
Set-TextValue $ws.Range("D13") "38.18"
$ws.Range("E13").Value = "  -17.82%  "

# Row 14
$ws.Range("E14").Value = "  -12.47%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.778.20"
$ws.Range("E15").Value = "  -7.63%  "

# Row 16
Set-TextValue $ws.Range("D16") "67.167.26"
$ws.Range("E16").Value = "  -5.01%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.260.10"
$ws.Range("E17").Value = "  -7.43%  "

# Row 18
Set-TextValue $ws.Range("D18") "535.74"
$ws.Range("E18").Value = "  -12.09%  "

# Row 19
$ws.Range("E19").Value = "  -6.25%  "

# Row 20
Set-TextValue $ws.Range("D20") "7.20"
$ws.Range("E20").Value = "  -14.14%  "

# Row 21
Set-TextValue $ws.Range("D21") "15.08"
$ws.Range("E21").Value = "  -14.96%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.760"
$ws.Range("E22").Value = "  -13.86%  "

# Row 23
Set-TextValue $ws.Range("D23") "7.88"
$ws.Range("E23").Value = "  -13.45%  "

# Row 24
Set-TextValue $ws.Range("D24") "85.18"
$ws.Range("E24").Value = "  -12.49%  "

# Row 25
Set-TextValue $ws.Range("D25") "13.54"
$ws.Range("E25").Value = "  -13.12%  "

# Row 26
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("E27").Value = "  -12.67%  "

# Row 28
Set-TextValue $ws.Range("D28") "29.31"
$ws.Range("E28").Value = "  -12.64%  "

# Row 29
$ws.Range("E29").Value = "  -11.64%  "

# Row 30
$ws.Range("E30").Value = "  -17.28%  "

# Row 31
$ws.Range("E31").Value = "  -10.94%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.14"
$ws.Range("E32").Value = "  -12.11%  "

# Row 33
$ws.Range("E33").Value = "  -17.83%  "

# Row 34
Set-TextValue $ws.Range("D34") "541.35"
$ws.Range("E34").Value = "  -14.99%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.69"
$ws.Range("E35").Value = "  -16.67%  "

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.0449"
$ws.Range("E37").Value = "  -7.67%  "

# Row 38
Set-TextValue $ws.Range("D38") "53.32"
$ws.Range("E38").Value = "  -5.96%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0854"
$ws.Range("E39").Value = "  -13.87%  "

# Row 40
Set-TextValue $ws.Range("D40") "9.09"
$ws.Range("E40").Value = "  -15.89%  "

# Row 41
$ws.Range("E41").Value = "  -10.21%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D42") "2.926.66"
$ws.Range("E42").Value = "  -12.57%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D43") "2.71"
$ws.Range("E43").Value = "  -20.96%  "

# Row 44
$ws.Range("E44").Value = "  -16.08%  "

# Row 45
$ws.Range("E45").Value = "  -19.20%  "

# Row 46
Set-TextValue $ws.Range("D46") "26.55"
$ws.Range("E46").Value = "  -16.97%  "

# Row 47
$ws.Range("E47").Value = "  -14.63%  "

# Row 48
$ws.Range("E48").Value = "  -0.09%  "

# Row 49
Set-TextValue $ws.Range("D49") "127.34"
$ws.Range("E49").Value = "  -4.25%  "

# Row 50
$ws.Range("E50").Value = "  -21.96%  "

# Row 51
$ws.Range("E51").Value = "  -12.79%  "
